$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1853085.1
$ws.Range("I28").Value = 2778127.8
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 2778127.8
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = -2777642.8
$ws.Range("N28").Value = -3970
$ws.Range("H33").Value = 77.708336
$ws.Range("I33").Value = 66.8421
$ws.Range("J33").Value = 119
$ws.Range("K33").Value = 66.8421
$ws.Range("L33").Value = 119
$ws.Range("M33").Value = 162.1579
$ws.Range("N33").Value = -577
$ws.Range("H38").Value = 75.166664
$ws.Range("I38").Value = 88.2
$ws.Range("J38").Value = 10
$ws.Range("K38").Value = 264.6
$ws.Range("L38").Value = 30
$ws.Range("M38").Value = 107.4
$ws.Range("N38").Value = -774
$ws.Range("H41").Value = 6944719.5
$ws.Range("J41").Value = 261
$ws.Range("L41").Value = 261
$ws.Range("N41").Value = -1141
$ws.Range("H103").Value = 605.9
$ws.Range("I103").Value = 572.875
$ws.Range("J103").Value = 617.9091
$ws.Range("K103").Value = 1718.625
$ws.Range("L103").Value = 1853.7273
$ws.Range("M103").Value = -1132.625
$ws.Range("N103").Value = -3025.7273
$ws.Range("H112").Value = 12398012
$ws.Range("J112").Value = 13637713
$ws.Range("L112").Value = 40913139
$ws.Range("N112").Value = -40915355
$ws.Range("H125").Value = 11211650
$ws.Range("I125").Value = 600
$ws.Range("J125").Value = 22422700
$ws.Range("K125").Value = 5400
$ws.Range("L125").Value = 201804300
$ws.Range("M125").Value = -2940
$ws.Range("N125").Value = -201809220
$ws.Range("H132").Value = 472162.22
$ws.Range("I132").Value = 579448.4
$ws.Range("J132").Value = 21560.4
$ws.Range("K132").Value = 1738345.2
$ws.Range("L132").Value = 64681.2
$ws.Range("M132").Value = -1735815.2
$ws.Range("N132").Value = -69741.20000000001
$ws.Range("H135").Value = 1464.4166
$ws.Range("I135").Value = 1464.4166
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13179.7494
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -10644.7494
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 35715336
$ws.Range("I137").Value = 52632250
$ws.Range("J137").Value = 1855.5555
$ws.Range("K137").Value = 157896750
$ws.Range("L137").Value = 5566.666499999999
$ws.Range("M137").Value = -157894200
$ws.Range("N137").Value = -10666.6665
$ws.Range("H138").Value = 5604766.5
$ws.Range("I138").Value = 1833638.9
$ws.Range("J138").Value = 7939274.5
$ws.Range("K138").Value = 5500916.699999999
$ws.Range("L138").Value = 23817823.5
$ws.Range("M138").Value = -5495776.699999999
$ws.Range("N138").Value = -23828103.5
$ws.Range("H141").Value = 3164
$ws.Range("I141").Value = 2299.889
$ws.Range("J141").Value = 7052.5
$ws.Range("K141").Value = 6899.667
$ws.Range("L141").Value = 21157.5
$ws.Range("M141").Value = -1719.667
$ws.Range("N141").Value = -31517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22471.346
$ws.Range("I32").Value = 3788.4314
$ws.Range("K32").Value = 3788.4314
$ws.Range("M32").Value = -3501.4314
$ws.Range("H45").Value = 1224.75
$ws.Range("I45").Value = 1224.75
$ws.Range("K45").Value = 1224.75
$ws.Range("M45").Value = -847.75
$ws.Range("H74").Value = 8073.1577
$ws.Range("I74").Value = 1189.6
$ws.Range("J74").Value = 15721.556
$ws.Range("K74").Value = 1189.6
$ws.Range("L74").Value = 15721.556
$ws.Range("M74").Value = -315.5999999999999
$ws.Range("N74").Value = -17469.556
$ws.Range("H77").Value = 8073.1577
$ws.Range("I77").Value = 1189.6
$ws.Range("J77").Value = 15721.556
$ws.Range("K77").Value = 5948
$ws.Range("L77").Value = 78607.78
$ws.Range("M77").Value = -1580
$ws.Range("N77").Value = -87343.78
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3998.5715
$ws.Range("I122").Value = 2709.4443
$ws.Range("K122").Value = 8128.3329
$ws.Range("M122").Value = -5678.3329
$ws.Range("H132").Value = 2283.7708
$ws.Range("I132").Value = 1751.5278
$ws.Range("K132").Value = 5254.5834
$ws.Range("M132").Value = -2724.5834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2574.0833
$ws.Range("I134").Value = 1937.6
$ws.Range("J134").Value = 5756.5
$ws.Range("K134").Value = 5812.799999999999
$ws.Range("L134").Value = 17269.5
$ws.Range("M134").Value = -3277.799999999999
$ws.Range("N134").Value = -22339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3778.7346
$ws.Range("I31").Value = 1417.619
$ws.Range("K31").Value = 1417.619
$ws.Range("M31").Value = -1122.619
$ws.Range("H34").Value = 3778.7346
$ws.Range("I34").Value = 1417.619
$ws.Range("K34").Value = 1417.619
$ws.Range("M34").Value = -1215.619
$ws.Range("H132").Value = 3508.348
$ws.Range("I132").Value = 3108.9412
$ws.Range("J132").Value = 4640
$ws.Range("K132").Value = 9326.8236
$ws.Range("L132").Value = 13920
$ws.Range("M132").Value = -6796.8236
$ws.Range("N132").Value = -18980
$ws.Range("H134").Value = 3573.1
$ws.Range("I134").Value = 1824
$ws.Range("J134").Value = 5004.1816
$ws.Range("K134").Value = 5472
$ws.Range("L134").Value = 15012.5448
$ws.Range("M134").Value = -2937
$ws.Range("N134").Value = -20082.5448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H107").Value = 445.89285
$ws.Range("I107").Value = 474.05264
$ws.Range("K107").Value = 1422.15792
$ws.Range("M107").Value = 497.8420799999999
$ws.Range("H123").Value = 2158.889
$ws.Range("I123").Value = 430
$ws.Range("J123").Value = 2375
$ws.Range("K123").Value = 1290
$ws.Range("L123").Value = 7125
$ws.Range("M123").Value = 1160
$ws.Range("N123").Value = -12025
$ws.Range("H131").Value = 5651021
$ws.Range("I131").Value = 229.8
$ws.Range("J131").Value = 6174242.5
$ws.Range("K131").Value = 689.4000000000001
$ws.Range("L131").Value = 18522727.5
$ws.Range("M131").Value = 4350.6
$ws.Range("N131").Value = -18532807.5
$ws.Range("H133").Value = 5647.4287
$ws.Range("I133").Value = 4375
$ws.Range("J133").Value = 7344
$ws.Range("K133").Value = 13125
$ws.Range("L133").Value = 22032
$ws.Range("M133").Value = -8065
$ws.Range("N133").Value = -32152
$ws.Range("H136").Value = 3619.2974
$ws.Range("I136").Value = 915
$ws.Range("J136").Value = 3947.0908
$ws.Range("K136").Value = 2745
$ws.Range("L136").Value = 11841.2724
$ws.Range("M136").Value = 2355
$ws.Range("N136").Value = -22041.2724
$ws.Range("H137").Value = 7218575.5
$ws.Range("J137").Value = 147719
$ws.Range("L137").Value = 443157
$ws.Range("N137").Value = -453357
$ws.Range("H138").Value = 802
$ws.Range("I138").Value = 802
$ws.Range("K138").Value = 2406
$ws.Range("M138").Value = 2734
$ws.Range("H139").Value = 1606.9
$ws.Range("I139").Value = 1606.9
$ws.Range("K139").Value = 4820.700000000001
$ws.Range("M139").Value = 319.2999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1391142.6
$ws.Range("I122").Value = 2778977.8
$ws.Range("K122").Value = 8336933.399999999
$ws.Range("M122").Value = -8334483.399999999
$ws.Range("H123").Value = 12443.714
$ws.Range("J123").Value = 12443.714
$ws.Range("L123").Value = 12443.714
$ws.Range("N123").Value = -17343.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2441.1428
$ws.Range("I68").Value = 2200
$ws.Range("J68").Value = 2762.6667
$ws.Range("K68").Value = 2200
$ws.Range("L68").Value = 2762.6667
$ws.Range("M68").Value = -1451
$ws.Range("N68").Value = -4260.6667
$ws.Range("H71").Value = 2441.1428
$ws.Range("I71").Value = 2200
$ws.Range("J71").Value = 2762.6667
$ws.Range("K71").Value = 11000
$ws.Range("L71").Value = 13813.3335
$ws.Range("M71").Value = -7256
$ws.Range("N71").Value = -21301.3335
$ws.Range("H140").Value = 44943
$ws.Range("J140").Value = 44943
$ws.Range("L140").Value = 44943
$ws.Range("N140").Value = -55303
$ws.Range("H141").Value = 45523
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 48153.75
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 48153.75
$ws.Range("N141").Value = -58513.75
$ws.Range("M141").Value = -29820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 19271
$ws.Range("J68").Value = 19271
$ws.Range("L68").Value = 19271
$ws.Range("N68").Value = -20893
$ws.Range("H71").Value = 19271
$ws.Range("J71").Value = 19271
$ws.Range("L71").Value = 57813
$ws.Range("N71").Value = -65925
$ws.Range("H100").Value = 495
$ws.Range("I100").Value = 300
$ws.Range("J100").Value = 560
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 1120
$ws.Range("M100").Value = -59
$ws.Range("N100").Value = -2202
$ws.Range("H132").Value = 5866.4443
$ws.Range("I132").Value = 6273.1333
$ws.Range("K132").Value = 18819.3999
$ws.Range("M132").Value = -16289.3999
$ws.Range("H136").Value = 2747.2942
$ws.Range("I136").Value = 1960.4
$ws.Range("K136").Value = 5881.200000000001
$ws.Range("M136").Value = -3331.200000000001
